# Fruta / hortaliza, semanal
# Insert two new weekly data rows (2021-10-18) after the current row 11,
# shifting the existing rows 12-41 down to 14-43.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before row 12 (existing data shifts down).
$ws.Rows.Item(12).Resize(2).Insert()

# New row 12: Región de O'Higgins, 2021-10-18
$ws.Range("A12").Value = 5
$ws.Range("B12").Value = "Macroferia Regional de Talca"
$ws.Range("C12").Value = "Maule"
$ws.Range("D12").Value = 44487
$ws.Range("E12").Value = 7
$ws.Range("F12").Value = 100112022
$ws.Range("G12").Value = "Arveja Verde"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 20000
$ws.Range("L12").Value = 20000
$ws.Range("M12").Value = 20000
$ws.Range("N12").Value = "`$/saco 25 kilos"
$ws.Range("O12").Value = "Región de O'Higgins"
$ws.Range("P12").Value = 800
$ws.Range("Q12").Value = 25
$ws.Range("R12").Value = "Hortaliza"

# New row 13: Región del Maule, 2021-10-18
$ws.Range("A13").Value = 5
$ws.Range("B13").Value = "Macroferia Regional de Talca"
$ws.Range("C13").Value = "Maule"
$ws.Range("D13").Value = 44487
$ws.Range("E13").Value = 7
$ws.Range("F13").Value = 100112022
$ws.Range("G13").Value = "Arveja Verde"
$ws.Range("H13").Value = "Sin especificar"
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 150
$ws.Range("K13").Value = 19000
$ws.Range("L13").Value = 19000
$ws.Range("M13").Value = 19000
$ws.Range("N13").Value = "`$/saco 25 kilos"
$ws.Range("O13").Value = "Región del Maule"
$ws.Range("P13").Value = 760
$ws.Range("Q13").Value = 25
$ws.Range("R13").Value = "Hortaliza"
